$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update B1 from a literal value to a formula (160*10 = 1600)
$ws.Range("B1").Formula = "=160*10"

# Move the active selection to B3 (cosmetic, matches author's cursor position)
$ws.Range("B3").Select()
